$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cells C2:C26 from Excel serial date 45253 (2023-11-23) to 45254 (2023-11-24)
for ($row = 2; $row -le 26; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45253) {
        $cell.Value2 = 45254
    }
}
